$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.961.79'
$ws.Range('E2').Value = '  -1.96%  '

$ws.Range('D3').Value = '2.171.09'
$ws.Range('E3').Value = '  -2.79%  '

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').Value = '''247.23'
$ws.Range('E5').Value = '  -2.09%  '

$ws.Range('E6').Value = '  -2.45%  '

$ws.Range('D7').Value = '''66.35'
$ws.Range('E7').Value = '  -7.34%  '

$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').Value = '''0.565'
$ws.Range('E9').Value = '  -0.17%  '

$ws.Range('D10').Value = '''57.94'
$ws.Range('E10').Value = '  -0.87%  '

$ws.Range('D11').Value = '''0.0925'
$ws.Range('E11').Value = '  -4.81%  '

$ws.Range('D12').Value = '''35.65'
$ws.Range('E12').Value = '  -15.39%  '

$ws.Range('D13').Value = '''0.104'
$ws.Range('E13').Value = '  -1.39%  '

$ws.Range('D14').Value = '''6.92'
$ws.Range('E14').Value = '  -0.65%  '

$ws.Range('E15').Value = '  -2.74%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '''0.854'
$ws.Range('E16').Value = '  -0.82%  '

$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '''14.23'
$ws.Range('E17').Value = '  -5.53%  '

$ws.Range('D18').Value = '2.181.18'
$ws.Range('E18').Value = '  -2.44%  '

$ws.Range('D19').Value = '40.863.40'
$ws.Range('E19').Value = '  -2.14%  '

$ws.Range('D20').Value = '0.0₃0936'
$ws.Range('E20').Value = '  -3.52%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''6.09'
$ws.Range('E21').Value = '  -1.82%  '

$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').Value = '''71.43'
$ws.Range('E22').Value = '  -2.41%  '

$ws.Range('D23').Value = '''229.54'
$ws.Range('E23').Value = '  -2.39%  '

$ws.Range('D24').Value = '''2.06'
$ws.Range('E24').Value = '  -8.73%  '

$ws.Range('D25').Value = '''11.72'
$ws.Range('E25').Value = '  +14.30%  '

$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.07%  '

$ws.Range('E27').Value = '  -0.52%  '

$ws.Range('D28').Value = '''2.40'
$ws.Range('E28').Value = '  -4.41%  '

$ws.Range('E29').Value = '  -3.49%  '

$ws.Range('D30').Value = '''168.34'
$ws.Range('E30').Value = '  -1.08%  '

$ws.Range('D31').Value = '''20.19'
$ws.Range('E31').Value = '  -2.89%  '

$ws.Range('E32').Value = '  -2.67%  '

$ws.Range('D33').Value = '''5.60'
$ws.Range('E33').Value = '  +2.44%  '

$ws.Range('D34').Value = '''0.0733'
$ws.Range('E34').Value = '  +1.55%  '

$ws.Range('E35').Value = '  -3.22%  '

$ws.Range('D36').Value = '''4.54'
$ws.Range('E36').Value = '  -2.78%  '

$ws.Range('D37').Value = '''25.52'
$ws.Range('E37').Value = '  -4.84%  '

$ws.Range('D38').Value = '''4.04'
$ws.Range('E38').Value = '  -2.15%  '

$ws.Range('D39').Value = '''0.0299'
$ws.Range('E39').Value = '  +5.51%  '

$ws.Range('E40').Value = '  -5.95%  '

$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = '''11.66'
$ws.Range('E41').Value = '  -0.25%  '

$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').Value = '''5.46'
$ws.Range('E42').Value = '  -9.60%  '

$ws.Range('D43').Value = '''60.06'
$ws.Range('E43').Value = '  -13.79%  '

$ws.Range('D44').Value = '''4.74'
$ws.Range('E44').Value = '  -6.94%  '

$ws.Range('E45').Value = '  -11.39%  '

$ws.Range('B46').Value = 'BinanceUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D46').Value = '''1.01'
$ws.Range('E46').Value = '  +0.47%  '

$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '''8.46'
$ws.Range('E47').Value = '  -4.93%  '

$ws.Range('D48').Value = '''0.0985'
$ws.Range('E48').Value = '  -3.48%  '

$ws.Range('E49').Value = '  -0.88%  '

$ws.Range('E50').Value = '  -3.70%  '

$ws.Range('E51').Value = '  -1.03%  '
